$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: add the "phase_2" value for S1 (new column F2)
$ws.Range("F2").Value = "10: 0 - 0 - 0 - 0"

# Row 3: S2 now only has a 2-phase layout (N-S), so update B3:E3 and clear F3
$ws.Range("B3").Value = "N - S"
$ws.Range("C3").Value = "20 - 20"
$ws.Range("D3").Value = "30 : 1 - 1"
$ws.Range("E3").Value = "40 : 0 - 0"
$ws.Range("F3").Clear()

# Update the view: zoom to 190% and select B2
$excel.ActiveWindow.Zoom = 190
$ws.Range("B2").Select()
